# Insert a new data row before the current row 346 (Macroferia Regional de
# Talca - Zanahoria), shifting the existing rows 346-424 down to 347-425.
# EntireRow.Insert() keeps formatting/styles on the shifted rows (e.g. the
# date style on column D), which is what we want since every later row in
# the diff is simply the previous row's data moved down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(346).Insert()

# Populate the newly-inserted row 346 with its own values.
$ws.Cells.Item(346, 1).Value  = 5
$ws.Cells.Item(346, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(346, 3).Value  = "Maule"
$ws.Cells.Item(346, 4).Value  = 44889
$ws.Cells.Item(346, 5).Value  = 7
$ws.Cells.Item(346, 6).Value  = 100114013
$ws.Cells.Item(346, 7).Value  = "Zanahoria"
$ws.Cells.Item(346, 8).Value  = "Sin especificar"
$ws.Cells.Item(346, 9).Value  = "Primera"
$ws.Cells.Item(346, 10).Value = 500
$ws.Cells.Item(346, 11).Value = 10000
$ws.Cells.Item(346, 12).Value = 10000
$ws.Cells.Item(346, 13).Value = 10000
$ws.Cells.Item(346, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(346, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(346, 16).Value = 500
$ws.Cells.Item(346, 17).Value = 20
$ws.Cells.Item(346, 18).Value = "Hortaliza"

# Make sure the new row's date cell (D346) carries the same number format
# as the rest of column D (matching the other rows' date style).
$ws.Cells.Item(346, 4).NumberFormat = $ws.Cells.Item(347, 4).NumberFormat()
